# Update field log for bottom depth
# Inserts a new "Bottom Depth" column (in meters) just before the existing
# "Comments" column, shifting Comments one column to the right, and tidies
# up the surrounding column widths / selection to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L (12); this shifts the old "Comments" column (and
# its merged header/banner range) one column to the right to M, and expands
# the A1:L3 banner merge + used range automatically to column M.
$ws.Columns.Item(12).Insert()

# Populate the header (row 5) and unit sub-label (row 6) for the new column.
$ws.Cells.Item(5, 12).Value = "Bottom Depth"
$ws.Cells.Item(6, 12).Value = "(m)"

# Match the final column widths: the "Start Station..." comment column (C)
# narrows slightly, the new Bottom Depth column is narrow, and the
# relocated Comments column (M) keeps most of the old width budget.
$ws.Columns.Item(3).ColumnWidth = 31.5
$ws.Columns.Item(12).ColumnWidth = 9.83203125
$ws.Columns.Item(13).ColumnWidth = 51.83203125

# Mirror the saved cursor/selection position.
$ws.Range("D8").Select()
